$wb = $excel.ActiveWorkbook

$newNames = @(
    "summ49401962",
    "summ49520588",
    "summ49622465",
    "summ49725367",
    "summ49822326",
    "summ49923252",
    "summ50022188",
    "summ50123029",
    "summ50225579",
    "summ50389253",
    "summ50500292",
    "summ50610560",
    "summ50718558",
    "summ50827392",
    "summ50930329",
    "summ51022285",
    "summ51132405",
    "summ51222819",
    "summ51323062",
    "summ51423363",
    "summ51522815",
    "summ51623435",
    "summ51722850",
    "summ51822947",
    "summ51928940",
    "summ52022608",
    "summ52122925",
    "summ52227908",
    "summ52326247",
    "summ52422732",
    "summ52521433",
    "summ52624906",
    "summ52721405",
    "summ52822337",
    "summ52924446",
    "summ53033892",
    "summ53131561",
    "summ53228943",
    "summ53340787",
    "summ53440607",
    "summ53547872",
    "summ53639103",
    "summ53744713",
    "summ53853899",
    "summ53957683",
    "summ54055873",
    "summ54164321",
    "summ54257928",
    "summ54360641",
    "summ54464956"
)

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Name = $newNames[$i - 1]
}

